$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2550667226314545
$ws.Range("B1").Value = 0.2074719965457916
$ws.Range("C1").Value = 0.1840784400701523
$ws.Range("D1").Value = 0.2016288042068481
$ws.Range("E1").Value = 0.2497632056474686
